# Generate Report for Archive
# - Update the localization status from "Ready for handoff" to "In Translation"
#   on the Overview sheet (zh-cn / de-de status columns) and on the per-locale
#   "zh-cn" / "de-de" sheets (Status column).
# - Narrow the status column(s) to match the shorter label.

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"

# --- Overview sheet: status columns E (zh-cn) and F (de-de) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

# --- zh-cn sheet: Status column C ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5

# --- de-de sheet: Status column C ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
